$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.05649126858622
$ws.Range("C2").Value = 11.11937511077658
$ws.Range("D2").Value = 14.9809489550595
$ws.Range("E2").Value = 16.39770022460419
$ws.Range("G2").Value = 3.667884810023161
$ws.Range("J2").Value = 9.361415295849442
$ws.Range("O2").Value = 27.33609796950253
$ws.Range("B3").Value = 17.45358153290566
$ws.Range("C3").Value = 10.55984671663852
$ws.Range("D3").Value = 14.92190621030298
$ws.Range("E3").Value = 16.33920153341992
$ws.Range("G3").Value = 3.670909957096306
$ws.Range("J3").Value = 9.371251528264512
$ws.Range("O3").Value = 27.43216004437329
$ws.Range("B4").Value = 17.07496380895604
$ws.Range("C4").Value = 10.20107496125499
$ws.Range("D4").Value = 14.8891232366978
$ws.Range("E4").Value = 16.30711106855635
$ws.Range("G4").Value = 3.672863441997519
$ws.Range("J4").Value = 9.378787643960573
$ws.Range("O4").Value = 27.49997969945234
$ws.Range("B5").Value = 16.91878789247725
$ws.Range("C5").Value = 10.0512012370221
$ws.Range("D5").Value = 14.87664525317529
$ws.Range("E5").Value = 16.29500479929133
$ws.Range("G5").Value = 3.673683738012624
$ws.Range("J5").Value = 9.382234851908871
$ws.Range("O5").Value = 27.52982527496845
$ws.Range("B6").Value = 16.89274829469804
$ws.Range("C6").Value = 10.02609816629172
$ws.Range("D6").Value = 14.87462676258906
$ws.Range("E6").Value = 16.29305341854729
$ws.Range("G6").Value = 3.673821413832912
$ws.Range("J6").Value = 9.382829972329398
$ws.Range("O6").Value = 27.53491408978265
$ws.Range("B7").Value = 17.07286487975987
$ws.Range("C7").Value = 10.1990683512066
$ws.Range("D7").Value = 14.88895137479593
$ws.Range("E7").Value = 16.30694385815608
$ws.Range("G7").Value = 3.672874406555283
$ws.Range("J7").Value = 9.378832611347953
$ws.Range("O7").Value = 27.50037328367503
$ws.Range("B8").Value = 17.85049116449869
$ws.Range("C8").Value = 10.92969369602329
$ws.Range("D8").Value = 14.95987660127699
$ws.Range("E8").Value = 16.37674070194101
$ws.Range("G8").Value = 3.66890800166384
$ws.Range("J8").Value = 9.364496127346147
$ws.Range("O8").Value = 27.3673779610746
$ws.Range("B9").Value = 19.29930329817671
$ws.Range("C9").Value = 12.23611796749721
$ws.Range("D9").Value = 15.12604789869196
$ws.Range("E9").Value = 16.54356895108348
$ws.Range("G9").Value = 3.661887896500161
$ws.Range("J9").Value = 9.348262612728591
$ws.Range("O9").Value = 27.17728077810643
$ws.Range("B10").Value = 20.3063855622829
$ws.Range("C10").Value = 13.11296542338122
$ws.Range("D10").Value = 15.26398500492358
$ws.Range("E10").Value = 16.68376900558137
$ws.Range("G10").Value = 3.657186759239508
$ws.Range("J10").Value = 9.343582601726053
$ws.Range("O10").Value = 27.08148785027336
$ws.Range("B11").Value = 20.75014474805262
$ws.Range("C11").Value = 13.49294072959859
$ws.Range("D11").Value = 15.3300160086012
$ws.Range("E11").Value = 16.75122117885762
$ws.Range("G11").Value = 3.655146034138953
$ws.Range("J11").Value = 9.343026408037721
$ws.Range("O11").Value = 27.0475864925979
$ws.Range("B12").Value = 20.91598152325679
$ws.Range("C12").Value = 13.63405138221135
$ws.Range("D12").Value = 15.35547716196362
$ws.Range("E12").Value = 16.77727711513413
$ws.Range("G12").Value = 3.654387243908881
$ws.Range("J12").Value = 9.343041714030152
$ws.Range("O12").Value = 27.03615149879456
$ws.Range("B13").Value = 20.88036572496851
$ws.Range("C13").Value = 13.60378503651683
$ws.Range("D13").Value = 15.34997359169703
$ws.Range("E13").Value = 16.77164292604504
$ws.Range("G13").Value = 3.654550042140067
$ws.Range("J13").Value = 9.343028374314533
$ws.Range("O13").Value = 27.03855166979942
$ws.Range("B14").Value = 20.76383317354769
$ws.Range("C14").Value = 13.50460595491924
$ws.Range("D14").Value = 15.33210165234739
$ws.Range("E14").Value = 16.75335461973556
$ws.Range("G14").Value = 3.655083328136413
$ws.Range("J14").Value = 9.343023141093369
$ws.Range("O14").Value = 27.04661755370529
$ws.Range("B15").Value = 20.69216259011004
$ws.Range("C15").Value = 13.44349260677243
$ws.Range("D15").Value = 15.32121357072374
$ws.Range("E15").Value = 16.74221888151931
$ws.Range("G15").Value = 3.65541180044089
$ws.Range("J15").Value = 9.343049349020163
$ws.Range("O15").Value = 27.05174114070707
$ws.Range("B16").Value = 20.27708346860096
$ws.Range("C16").Value = 13.08774737806175
$ws.Range("D16").Value = 15.25973448261776
$ws.Range("E16").Value = 16.67943351562056
$ws.Range("G16").Value = 3.657322087310695
$ws.Range("J16").Value = 9.343650580166026
$ws.Range("O16").Value = 27.08389906952166
$ws.Range("B17").Value = 20.01866166464467
$ws.Range("C17").Value = 12.86462091655739
$ws.Range("D17").Value = 15.22284914670651
$ws.Range("E17").Value = 16.64184738442097
$ws.Range("G17").Value = 3.65851898859226
$ws.Range("J17").Value = 9.344422099278631
$ws.Range("O17").Value = 27.10611336693398
$ws.Range("B18").Value = 19.86867915236203
$ws.Range("C18").Value = 12.73450756320723
$ws.Range("D18").Value = 15.20194364075384
$ws.Range("E18").Value = 16.62057566044544
$ws.Range("G18").Value = 3.659216629145233
$ws.Range("J18").Value = 9.345013912932977
$ws.Range("O18").Value = 27.11980049615781
$ws.Range("B19").Value = 19.8176711374216
$ws.Range("C19").Value = 12.69015022169142
$ws.Range("D19").Value = 15.19491908769156
$ws.Range("E19").Value = 16.61343343310245
$ws.Range("G19").Value = 3.659454423368976
$ws.Range("J19").Value = 9.345239725400997
$ws.Range("O19").Value = 27.12459066478402
$ws.Range("B20").Value = 20.04631138282046
$ws.Range("C20").Value = 12.88855747715686
$ws.Range("D20").Value = 15.22674368257172
$ws.Range("E20").Value = 16.64581270384538
$ws.Range("G20").Value = 3.65839062333052
$ws.Range("J20").Value = 9.344324647943877
$ws.Range("O20").Value = 27.10365434516268
$ws.Range("B21").Value = 20.79812248477004
$ws.Range("C21").Value = 13.53381305540549
$ws.Range("D21").Value = 15.33733880781989
$ws.Range("E21").Value = 16.75871253386343
$ws.Range("G21").Value = 3.654926310131374
$ws.Range("J21").Value = 9.343018549038508
$ws.Range("O21").Value = 27.04421025173415
$ws.Range("B22").Value = 21.27656769411631
$ws.Range("C22").Value = 13.93931443433034
$ws.Range("D22").Value = 15.41227308030561
$ws.Range("E22").Value = 16.83548314211866
$ws.Range("G22").Value = 3.652743676624548
$ws.Range("J22").Value = 9.343481623663004
$ws.Range("O22").Value = 27.01354092584343
$ws.Range("B23").Value = 21.02243524413979
$ws.Range("C23").Value = 13.72439055453145
$ws.Range("D23").Value = 15.37204171062692
$ws.Range("E23").Value = 16.79424141124039
$ws.Range("G23").Value = 3.65390115936831
$ws.Range("J23").Value = 9.343114097331311
$ws.Range("O23").Value = 27.02915761045172
$ws.Range("B24").Value = 20.03381533372024
$ws.Range("C24").Value = 12.87774146912008
$ws.Range("D24").Value = 15.22498202616889
$ws.Range("E24").Value = 16.64401893197521
$ws.Range("G24").Value = 3.658448627586489
$ws.Range("J24").Value = 9.344368243877256
$ws.Range("O24").Value = 27.10476321701807
$ws.Range("B25").Value = 18.91671680377606
$ws.Range("C25").Value = 11.89690941904891
$ws.Range("D25").Value = 15.0782585026765
$ws.Range("E25").Value = 16.49529115794959
$ws.Range("G25").Value = 3.663706448320891
$ws.Range("J25").Value = 9.351381572836166
$ws.Range("O25").Value = 27.22105585350922
